$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting the existing rows 3 and 4 down to 4 and 5.
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the CSKA Sofia vs Botev Vratsa match data.
$ws.Cells.Item(3,1).Value = "baEjVDx0"  # A3
$ws.Cells.Item(3,2).Value = "24/10/2024"  # B3
$ws.Cells.Item(3,3).Value = "11:30"  # C3
$ws.Cells.Item(3,4).Value = "BULGARIA - PARVA LIGA"  # D3
$ws.Cells.Item(3,5).Value = "CSKA Sofia"  # E3
$ws.Cells.Item(3,6).Value = "Botev Vratsa"  # F3
$ws.Cells.Item(3,7).Value = 1.27  # G3
$ws.Cells.Item(3,8).Value = 5.5  # H3
$ws.Cells.Item(3,9).Value = 10  # I3
$ws.Cells.Item(3,10).Value = 1.73  # J3
$ws.Cells.Item(3,11).Value = 2.4  # K3
$ws.Cells.Item(3,12).Value = 10  # L3
$ws.Cells.Item(3,13).Value = 1.06  # M3
$ws.Cells.Item(3,14).Value = 10  # N3
$ws.Cells.Item(3,15).Value = 1.25  # O3
$ws.Cells.Item(3,16).Value = 3.75  # P3
$ws.Cells.Item(3,17).Value = 1.88  # Q3
$ws.Cells.Item(3,18).Value = 1.98  # R3
$ws.Cells.Item(3,19).Value = 1.36  # S3
$ws.Cells.Item(3,20).Value = 3  # T3
$ws.Cells.Item(3,21).Value = 2.5  # U3
$ws.Cells.Item(3,22).Value = 1.5  # V3
$ws.Cells.Item(3,23).Value = 6  # W3
$ws.Cells.Item(3,24).Value = 5.5  # X3
$ws.Cells.Item(3,25).Value = 9.5  # Y3
$ws.Cells.Item(3,26).Value = 7  # Z3
$ws.Cells.Item(3,27).Value = 13  # AA3
$ws.Cells.Item(3,28).Value = 41  # AB3
$ws.Cells.Item(3,29).Value = 10  # AC3
$ws.Cells.Item(3,30).Value = 11  # AD3
$ws.Cells.Item(3,31).Value = 29  # AE3
$ws.Cells.Item(3,32).Value = 101  # AF3
$ws.Cells.Item(3,33).Value = 201  # AG3
$ws.Cells.Item(3,34).Value = 19  # AH3
$ws.Cells.Item(3,35).Value = 51  # AI3
$ws.Cells.Item(3,36).Value = 29  # AJ3
$ws.Cells.Item(3,37).Value = 151  # AK3
$ws.Cells.Item(3,38).Value = 81  # AL3
$ws.Cells.Item(3,39).Value = 81  # AM3
$ws.Cells.Item(3,40).Value = 3  # AN3
$ws.Cells.Item(3,41).Value = 6  # AO3
$ws.Cells.Item(3,42).Value = 21  # AP3
$ws.Cells.Item(3,43).Value = 17  # AQ3
$ws.Cells.Item(3,44).Value = 41  # AR3
$ws.Cells.Item(3,45).Value = 201  # AS3
$ws.Cells.Item(3,46).Value = 3  # AT3
$ws.Cells.Item(3,47).Value = 11  # AU3
$ws.Cells.Item(3,48).Value = 81  # AV3
$ws.Cells.Item(3,49).Value = 11  # AW3
$ws.Cells.Item(3,50).Value = 51  # AX3
$ws.Cells.Item(3,51).Value = 51  # AY3
$ws.Cells.Item(3,52).Value = 301  # AZ3
$ws.Cells.Item(3,53).Value = 351  # BA3
$ws.Cells.Item(3,54).Value = 501  # BB3
$ws.Cells.Item(3,55).Value = 51  # BC3
$ws.Cells.Item(3,56).Value = 51  # BD3
